$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (shifts existing rows 9-18 down to 10-19)
$ws.Rows.Item(9).Insert()

# Fill in new row 9 with the new PLAYER_DASH_COOLTIME parameter
$ws.Range("A9").Value = "PLAYER_DASH_COOLTIME"
$ws.Range("B9").Value = 30
$ws.Range("C9").Value = "int"

# Update selection to match target state
$ws.Range("K19").Select()
